{"js": "// Add a new paragraph style \"CompactList\" (display name \"Compact List\"),\n// mirroring the existing \"Compact\" style: based on \"Body Text\", marked as a\n// quick style, with 1.8pt (36 twentieths-of-a-point) spacing before/after.\n\n// Word's Styles collection always appends new styles at the end, so create\n// it first (name becomes \"Compact List\" -> styleId \"CompactList\").\ncontext.document.addStyle(\"Compact List\", Word.StyleType.paragraph);\nawait context.sync();\n\n// Re-fetch by name to get a stable, correctly-anchored reference before\n// mutating its properties (the object returned directly from addStyle can\n// lose track of which style it points to once other styles are touched).\nconst style = context.document.getStyles().getByName(\"Compact List\");\nstyle.baseStyle = \"BodyText\";\nstyle.quickStyle = true;\nstyle.paragraphFormat.spaceBefore = 1.8;\nstyle.paragraphFormat.spaceAfter = 1.8;\nawait context.sync();\n", "ps1": "# Add a new paragraph style \"CompactList\" (display name \"Compact List\"),\n# mirroring the existing \"Compact\" style: based on \"Body Text\", marked as a\n# quick style, with 1.8pt (36 twentieths-of-a-point) spacing before/after.\n\n$d = $word.ActiveDocument\n\n# wdStyleTypeParagraph = 1. Word's Styles collection always appends new\n# styles at the end (name \"Compact List\" -> styleId \"CompactList\").\n$s = $d.Styles.Add(\"Compact List\", 1)\n$s.BaseStyle = \"BodyText\"\n$s.QuickStyle = $true\n$s.ParagraphFormat.SpaceBefore = 1.8\n$s.ParagraphFormat.SpaceAfter = 1.8\n"}
